$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price column D, volume(1h) column E)
# D-column cells are forced to Text format so numeric-looking values
# (e.g. "9.30", "260.26") are preserved exactly as strings, matching source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.091.30"
$ws.Range("E2").Value = "  +3.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.217.21"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.26"
$ws.Range("E5").Value = "  +2.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "82.56"
$ws.Range("E6").Value = "  +11.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +2.79%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  +3.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.81"
$ws.Range("E10").Value = "  +6.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  +2.02%  "

$ws.Range("E12").Value = "  +3.89%  "

$ws.Range("E13").Value = "  +2.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.553.14"
$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.55"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.216.36"
$ws.Range("E16").Value = "  +2.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.783"
$ws.Range("E17").Value = "  +1.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.947.62"
$ws.Range("E18").Value = "  +4.02%  "

$ws.Range("E19").Value = "  +1.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.34"
$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("E21").Value = "  +1.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.37"
$ws.Range("E22").Value = "  +9.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.06"
$ws.Range("E23").Value = "  +2.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.30"
$ws.Range("E24").Value = "  -2.71%  "

$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.79"
$ws.Range("E26").Value = "  +2.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.23"
$ws.Range("E27").Value = "  +11.01%  "

$ws.Range("E28").Value = "  +1.13%  "

$ws.Range("E29").Value = "  +2.48%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").Value = "  +2.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.64"
$ws.Range("E32").Value = "  +2.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0876"
$ws.Range("E33").Value = "  +8.44%  "

$ws.Range("E34").Value = "  +3.92%  "

$ws.Range("E35").Value = "  +7.68%  "

$ws.Range("E36").Value = "  +1.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0361"
$ws.Range("E37").Value = "  +7.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.51"
$ws.Range("E38").Value = "  +6.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.59"
$ws.Range("E39").Value = "  +12.90%  "

$ws.Range("E40").Value = "  +18.00%  "

$ws.Range("E41").Value = "  +2.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.57"
$ws.Range("E42").Value = "  +8.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.14"
$ws.Range("E43").Value = "  +5.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.201"
$ws.Range("E44").Value = "  +2.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.99"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0989"
$ws.Range("E46").Value = "  +1.86%  "

$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.56"
$ws.Range("E48").Value = "  +29.25%  "

$ws.Range("E49").Value = "  +2.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.18"
$ws.Range("E50").Value = "  +4.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.443"
$ws.Range("E51").Value = "  -5.79%  "
